$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Second paragraph: swap the old instructional sentence for "For students:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This is the document you will need to change.  Delete everything below the above line.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For students:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Third paragraph ("Add instructions for your tutor to pull ...") becomes
#    four separate paragraphs. Do this as a single Find/Replace using "^p"
#    paragraph breaks in the replacement text (this keeps everything as
#    plain runs, splitting the original single paragraph into four, and
#    does not disturb the surrounding paragraphs).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Add instructions for your tutor to pull (merge) your request to the mainline.  Note that the tutor will not merge all requests, but you should say how it could be achieved.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To pull a project to the mainline, the project must be forked into the private repository, then the file must be edited and saved the individuals hard drive on their personal pc, then the file must be uploaded again to the private repository and then a pull request must be made to merge the document into the main repository.^p" +
    "For tutor:^p" +
    "For the tutor to accept and merge the pull request they must access the pull requests tab on the main screen of their major repository and accept the pull request from the user. (Tutor can acee and view changes in files, and review changes ect.)",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) The final paragraph needs finer-grained run structure: a proofErr
#    spell-check pair wrapped around "acee", and a bookmark ("_GoBack")
#    immediately before "ect.)". Replace that whole paragraph's content
#    (but not its trailing paragraph mark) via a flat-OPC WordOpenXML
#    fragment so the runs/bookmark/proofErr come out exactly right.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">For the tutor to accept and merge the pull request they must access the pull requests tab on the main screen of their major repository and accept the pull request from the user. (Tutor can </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>acee</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> and view changes in files, and review changes </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>ect.)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$finalRange.InsertXML($xml)
